# Added login test cases: rename sheet, add username/password + admin/pointofsale
# sample rows, resize the password column, and select/zoom on the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "ValidLogin"

# Write row 2 (credentials) before row 1 (headers) so the shared-string table
# is built in the same order as the source workbook: admin, pointofsale,
# username, password.
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "pointofsale"
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "password"

# Column B ("pointofsale") needs to be a bit wider to fit its contents.
$ws.Columns.Item(2).ColumnWidth = 11.140625

# Leave the selection on B2 and zoom the sheet in, matching the saved view.
$ws.Range("B2").Select()
$excel.ActiveWindow.Zoom = 190
